$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, pushing existing rows 118:141 down to 119:142
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new weekly record
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = 44511
$ws.Range("D118").Style = $ws.Range("D119").Style
$ws.Range("D118").NumberFormat = $ws.Range("D119").NumberFormat
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = 100112032
$ws.Range("G118").Value = "Zapallo italiano"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 100
$ws.Range("K118").Value = 8000
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 8500
$ws.Range("N118").Value = "$/caja 60 unidades"
$ws.Range("O118").Value = "Región del Maule"
$ws.Range("P118").Value = 142
$ws.Range("Q118").Value = 60
$ws.Range("R118").Value = "Hortaliza"
